$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header row: a new "Wheat.Ear.Wt" column is inserted right after
#    "Wheat.AboveGround.Wt" (column L). Grain.Wt / Grain.Size / Grain.Number
#    shift one column to the right (M->N, N->O, O->P).
# ---------------------------------------------------------------------------
$ws.Cells.Item(1, 13).Value() = "Wheat.Ear.Wt"       # M1
$ws.Cells.Item(1, 14).Value() = "Wheat.Grain.Wt"     # N1
$ws.Cells.Item(1, 15).Value() = "Wheat.Grain.Size"   # O1
$ws.Cells.Item(1, 16).Value() = "Wheat.Grain.Number" # P1

# ---------------------------------------------------------------------------
# 2. A handful of re-calculated Wheat.Phenology.HaunStage (column D) values.
# ---------------------------------------------------------------------------
$ws.Cells.Item(116, 4).Value() = 13.875
$ws.Cells.Item(163, 4).Value() = 10.66666666666667
$ws.Cells.Item(172, 4).Value() = 11.66666666666667
$ws.Cells.Item(181, 4).Value() = 11.58333333333333

# ---------------------------------------------------------------------------
# 3. Rows 206-213: the Ear.Wt observation that used to live in column K
#    (Spike.Wt) is copied into the new column M, leaving K untouched.
# ---------------------------------------------------------------------------
for ($r = 206; $r -le 213; $r++) {
    $ws.Cells.Item($r, 13).Value() = $ws.Cells.Item($r, 11).Value()  # M = K
}

# ---------------------------------------------------------------------------
# 4. Rows 214-221: same Ear.Wt observation, but this time it is moved out of
#    column K (K is cleared) into the new column M.
# ---------------------------------------------------------------------------
for ($r = 214; $r -le 221; $r++) {
    $ws.Cells.Item($r, 13).Value() = $ws.Cells.Item($r, 11).Value()  # M = K
    $ws.Cells.Item($r, 11).Value() = ""                              # clear K
}

# ---------------------------------------------------------------------------
# 5. Rows 222-229: columns M/N/O already held Grain.Wt/Grain.Size/
#    Grain.Number, and P held Ear.Wt. Shift them right by one column so the
#    Ear.Wt value (old P) lands in the new M column.
# ---------------------------------------------------------------------------
for ($r = 222; $r -le 229; $r++) {
    $oldM = $ws.Cells.Item($r, 13).Value()
    $oldN = $ws.Cells.Item($r, 14).Value()
    $oldO = $ws.Cells.Item($r, 15).Value()
    $oldP = $ws.Cells.Item($r, 16).Value()

    $ws.Cells.Item($r, 13).Value() = $oldP   # M = old P (Ear.Wt)
    $ws.Cells.Item($r, 14).Value() = $oldM   # N = old M (Grain.Wt)
    $ws.Cells.Item($r, 15).Value() = $oldN   # O = old N (Grain.Size)
    $ws.Cells.Item($r, 16).Value() = $oldO   # P = old O (Grain.Number)
}
